$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Gvu1234"
$ws.Range("A3").Value = "Gvu1235"
$ws.Range("A4").Value = "Gvu1236"
$ws.Range("A5").Value = "Gvu1237"
$ws.Range("A6").Value = "Gvu1238"

$ws.Range("A1:E6").Select()
